$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update the login credentials used by the test data sheet
$ws.Range("A2").Value = "swati.chetty@thomsonreuters.com.lrp5qa"
$ws.Range("B2").Value = "March12345"

# Turn the username cell into a mailto hyperlink (as Excel does automatically
# when a value looks like an email address)
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:swati.chetty@thomsonreuters.com.lrp5qa")

# Leave the active selection on the password cell
$ws.Range("B2").Select()
